# NIT-9012095521.xlsx — "Elimina EC anteriores y se agregan nuevos,
# se modifica base de datos"
#
# The period/value pairs in the detail table (rows 16-19) are re-listed in
# reverse period order, and the "Valor Mora" figures for the first/last
# rows are swapped accordingly:
#   before: 2304/42000  2305/60000  2306/60000  2307/22000
#   after : 2307/22000  2306/60000  2305/60000  2304/42000

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Periodo Mora column (E16:E19) -----------------------------------
$ws.Range("E16").Value = "2307"
$ws.Range("E17").Value = "2306"
$ws.Range("E18").Value = "2305"
$ws.Range("E19").Value = "2304"

# --- Valor Mora column (F16 / F19 swap) -------------------------------
$ws.Range("F16").Value = 22000
$ws.Range("F19").Value = 42000

# --- Remove the accent fill from the detail table header/body ---------
# (previously filled with a light theme accent colour; the refreshed
# template drops that shading so the rows sit on a plain white fill)
$ws.Range("B16:G19").Interior.Pattern = -4142   # xlPatternNone

# --- Column width refresh (auto-fit drift from the template update) ---
$ws.Columns.Item(2).ColumnWidth = 17.666666666666668   # B
$ws.Columns.Item(3).ColumnWidth = 15.833333333333334   # C
$ws.Columns.Item(5).ColumnWidth = 12.666666666666666   # E
$ws.Columns.Item(6).ColumnWidth = 9.333333333333332    # F
$ws.Columns.Item(7).ColumnWidth = 13.5                 # G
$ws.Columns.Item(8).ColumnWidth = 18.5                 # H
$ws.Columns.Item(9).ColumnWidth = 17.333333333333336   # I
$ws.Columns.Item(10).ColumnWidth = 14.166666666666666  # J
